$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.160.27"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.074.04"
$ws.Range("E3").Value = "  -1.10%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "253.91"
$ws.Range("E5").Value = "  +1.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.678"
$ws.Range("E6").Value = "  +2.62%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "62.16"
$ws.Range("E7").Value = "  +21.83%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.392"
$ws.Range("E9").Value = "  +4.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "61.63"
$ws.Range("E10").Value = "  +0.66%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0803"
$ws.Range("E11").Value = "  +8.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "16.39"
$ws.Range("E13").Value = "  +7.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.374.68"
$ws.Range("E14").Value = "  -1.17%  "
$ws.Range("E15").Value = "  -0.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.57"
$ws.Range("E16").Value = "  +9.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.068.85"
$ws.Range("E17").Value = "  -1.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.142.50"
$ws.Range("E18").Value = "  -0.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "15.67"
$ws.Range("E19").Value = "  +15.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "74.90"
$ws.Range("E20").Value = "  +4.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0929"
$ws.Range("E21").Value = "  +11.66%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.51"
$ws.Range("E22").Value = "  +6.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "241.28"
$ws.Range("E23").Value = "  +0.61%  "
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("E25").Value = "  -1.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.40"
$ws.Range("E26").Value = "  +20.53%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "171.16"
$ws.Range("E27").Value = "  +0.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.39"
$ws.Range("E28").Value = "  +2.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.48"
$ws.Range("E29").Value = "  -0.63%  "
$ws.Range("E30").Value = "  +3.49%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.88"
$ws.Range("E31").Value = "  +8.87%  "
$ws.Range("E32").Value = "  +4.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0642"
$ws.Range("E33").Value = "  +5.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.48"
$ws.Range("E34").Value = "  +9.72%  "
$ws.Range("E35").Value = "  -2.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.30"
$ws.Range("E37").Value = "  -0.77%  "
$ws.Range("B38").Value = "Cronos"
$ws.Range("C38").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.117"
$ws.Range("E38").Value = "  +28.98%  "
$ws.Range("B39").Value = "WEMIXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.77"
$ws.Range("E39").Value = "  -3.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.37"
$ws.Range("E40").Value = "  +4.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.27"
$ws.Range("E41").Value = "  +2.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0228"
$ws.Range("E42").Value = "  +2.13%  "
$ws.Range("E43").Value = "  +1.01%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.66"
$ws.Range("E44").Value = "  +1.07%  "
$ws.Range("B45").Value = "FTXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.50"
$ws.Range("E45").Value = "  +28.28%  "
$ws.Range("E46").Value = "  +2.94%  "
$ws.Range("E47").Value = "  +12.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.50"
$ws.Range("E48").Value = "  +9.43%  "
$ws.Range("B49").Value = "MXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.96"
$ws.Range("E49").Value = "  -2.12%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.308.33"
$ws.Range("E50").Value = "  -0.82%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.97"
$ws.Range("E51").Value = "  +0.06%  "
